# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for Kiwi (Vega Modelo de Temuco) above the
# existing data block (old row 434), shifting the existing rows 434-455 down
# to 438-459, then populate the newly inserted rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at row 434 (existing rows 434:455 shift down to 438:459)
$ws.Rows("434:437").Insert()

# Row 434: Kiwi Hayward - Especial
$ws.Range("A434").Value2 = 10
$ws.Range("B434").Value2 = "Vega Modelo de Temuco"
$ws.Range("C434").Value2 = "La Araucanía"
$ws.Range("D434").Value2 = 44746
$ws.Range("E434").Value2 = 9
$ws.Range("F434").Value2 = "Fruta"
$ws.Range("G434").Value2 = 100101
$ws.Range("H434").Value2 = "Berries"
$ws.Range("I434").Value2 = 100101007
$ws.Range("J434").Value2 = "Kiwi"
$ws.Range("K434").Value2 = "Hayward"
$ws.Range("L434").Value2 = "Especial"
$ws.Range("M434").Value2 = 600
$ws.Range("N434").Value2 = 10000
$ws.Range("O434").Value2 = 12000
$ws.Range("P434").Value2 = 11000
$ws.Range("Q434").Value2 = "$/bandeja 10 kilos"
$ws.Range("R434").Value2 = "Región de O'Higgins"
$ws.Range("S434").Value2 = 1100
$ws.Range("T434").Value2 = 10

# Row 435: Kiwi Hayward - Primera
$ws.Range("A435").Value2 = 10
$ws.Range("B435").Value2 = "Vega Modelo de Temuco"
$ws.Range("C435").Value2 = "La Araucanía"
$ws.Range("D435").Value2 = 44746
$ws.Range("E435").Value2 = 9
$ws.Range("F435").Value2 = "Fruta"
$ws.Range("G435").Value2 = 100101
$ws.Range("H435").Value2 = "Berries"
$ws.Range("I435").Value2 = 100101007
$ws.Range("J435").Value2 = "Kiwi"
$ws.Range("K435").Value2 = "Hayward"
$ws.Range("L435").Value2 = "Primera"
$ws.Range("M435").Value2 = 250
$ws.Range("N435").Value2 = 12000
$ws.Range("O435").Value2 = 12000
$ws.Range("P435").Value2 = 12000
$ws.Range("Q435").Value2 = "$/bandeja 18 kilos"
$ws.Range("R435").Value2 = "Región de O'Higgins"
$ws.Range("S435").Value2 = 667
$ws.Range("T435").Value2 = 18

# Row 436: Kiwi Hayward - Segunda
$ws.Range("A436").Value2 = 10
$ws.Range("B436").Value2 = "Vega Modelo de Temuco"
$ws.Range("C436").Value2 = "La Araucanía"
$ws.Range("D436").Value2 = 44746
$ws.Range("E436").Value2 = 9
$ws.Range("F436").Value2 = "Fruta"
$ws.Range("G436").Value2 = 100101
$ws.Range("H436").Value2 = "Berries"
$ws.Range("I436").Value2 = 100101007
$ws.Range("J436").Value2 = "Kiwi"
$ws.Range("K436").Value2 = "Hayward"
$ws.Range("L436").Value2 = "Segunda"
$ws.Range("M436").Value2 = 100
$ws.Range("N436").Value2 = 10000
$ws.Range("O436").Value2 = 10000
$ws.Range("P436").Value2 = 10000
$ws.Range("Q436").Value2 = "$/bandeja 18 kilos"
$ws.Range("R436").Value2 = "Región de O'Higgins"
$ws.Range("S436").Value2 = 556
$ws.Range("T436").Value2 = 18

# Row 437: Kiwi Hayward - Segunda (bins)
$ws.Range("A437").Value2 = 10
$ws.Range("B437").Value2 = "Vega Modelo de Temuco"
$ws.Range("C437").Value2 = "La Araucanía"
$ws.Range("D437").Value2 = 44746
$ws.Range("E437").Value2 = 9
$ws.Range("F437").Value2 = "Fruta"
$ws.Range("G437").Value2 = 100101
$ws.Range("H437").Value2 = "Berries"
$ws.Range("I437").Value2 = 100101007
$ws.Range("J437").Value2 = "Kiwi"
$ws.Range("K437").Value2 = "Hayward"
$ws.Range("L437").Value2 = "Segunda"
$ws.Range("M437").Value2 = 4
$ws.Range("N437").Value2 = 200000
$ws.Range("O437").Value2 = 200000
$ws.Range("P437").Value2 = 200000
$ws.Range("Q437").Value2 = "$/bins (450 kilos)"
$ws.Range("R437").Value2 = "Región de O'Higgins"
$ws.Range("S437").Value2 = 444
$ws.Range("T437").Value2 = 450
